$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap name/email values in A2/B2, also fixing the leading space on the name
$ws.Range("A2").Value = "shubham kumar"
$ws.Range("B2").Value = "shubhamk@gmail.com"

# Update the saved selection to A2
$ws.Activate()
[void]$ws.Range("A2").Select()
